$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H62" = 9367.75
    "I62" = 9367.75
    "K62" = 9367.75
    "M62" = -8743.75
    "H65" = 9367.75
    "I65" = 9367.75
    "K65" = 46838.75
    "M65" = -43718.75
    "H70" = 1669.1538
    "I70" = 733
    "J70" = 1950
    "K70" = 2199
    "L70" = 5850
    "M70" = -1929
    "N70" = -6390
    "H73" = 1669.1538
    "I73" = 733
    "J73" = 1950
    "K73" = 2199
    "L73" = 5850
    "M73" = -1263
    "N73" = -7722
    "H74" = 15630
    "I74" = 15630
    "K74" = 15630
    "M74" = -14694
    "H77" = 15630
    "I77" = 15630
    "K77" = 78150
    "M77" = -73470
    "H80" = 3119.2856
    "J80" = 3582.5
    "L80" = 10747.5
    "N80" = -12743.5
    "H83" = 3119.2856
    "J83" = 3582.5
    "L83" = 32242.5
    "N83" = -42226.5
    "H86" = 2759.625
    "I86" = 2510
    "J86" = 3009.25
    "K86" = 2510
    "L86" = 3009.25
    "M86" = -1387
    "N86" = -5255.25
    "H89" = 2759.625
    "I89" = 2510
    "J89" = 3009.25
    "K89" = 12550
    "L89" = 15046.25
    "M89" = -6934
    "N89" = -26278.25
    "H103" = 1472.1818
    "I103" = 1132
    "J103" = 1599.75
    "K103" = 3396
    "L103" = 4799.25
    "M103" = -2810
    "N103" = -5971.25
    "H112" = 2296.6924
    "I112" = 1285.6666
    "J112" = 2600
    "K112" = 3856.9998
    "L112" = 7800
    "M112" = -2748.9998
    "N112" = -10016
    "H113" = 4650
    "J113" = 4650
    "L113" = 4650
    "N113" = -11158
    "H115" = 2783.3333
    "I115" = 2499.2
    "J115" = 4204
    "K115" = 7497.599999999999
    "L115" = 12612
    "M115" = -5930.599999999999
    "H138" = 4049.6047
    "J138" = 4869.0967
    "L138" = 14607.2901
    "N138" = -24887.2901
    "N115" = -15746
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H31" = 21944
    "I31" = 12592
    "K31" = 12592
    "M31" = -12298
    "H32" = 8951.214
    "I32" = 8687.205
    "K32" = 8687.205
    "M32" = -8400.205
    "H45" = 2598
    "I45" = 2683.5715
    "J45" = 1999
    "K45" = 2683.5715
    "L45" = 1999
    "M45" = -2306.5715
    "N45" = -2753
    "H61" = 2571.0833
    "I61" = 2270.182
    "K61" = 2270.182
    "M61" = -2058.182
    "H102" = 1446.6
    "I102" = 1446.6
    "K102" = 1446.6
    "M102" = 175.4000000000001
    "H122" = 1746.6154
    "I122" = 1306.6364
    "K122" = 3919.9092
    "M122" = -1469.9092
    "H136" = 2571.0833
    "I136" = 2270.182
    "K136" = 6810.545999999999
    "M136" = -4260.545999999999
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H33" = 80000
    "I33" = 0
    "K33" = 0
    "H38" = 40000
    "J38" = 40000
    "L38" = 40000
    "N38" = -40832
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
$clears = @("M33")
foreach ($ref in $clears) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H38" = 29999.5
    "I38" = 19999
    "K38" = 19999
    "M38" = -19622
    "H46" = 29999.5
    "I46" = 19999
    "K46" = 19999
    "M46" = -19788
    "H62" = 5162
    "I62" = 3909.8333
    "K62" = 3909.8333
    "M62" = -3285.8333
    "H65" = 5162
    "I65" = 3909.8333
    "K65" = 19549.1665
    "M65" = -16429.1665
    "H141" = 236659.33
    "J141" = 236659.33
    "L141" = 236659.33
    "N141" = -247019.33
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H136" = 43555
    "J136" = 43555
    "L136" = 130665
    "N136" = -135765
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H7" = 3251
    "I7" = 3251
    "K7" = 3251
    "M7" = -3139
    "H68" = 2840.5
    "I68" = 2840.5
    "J68" = 0
    "K68" = 2840.5
    "L68" = 0
    "M68" = -2091.5
    "H71" = 2840.5
    "I71" = 2840.5
    "J71" = 0
    "K71" = 14202.5
    "L71" = 0
    "M71" = -10458.5
    "H106" = 63079.4
    "J106" = 63079.4
    "L106" = 63079.4
    "N106" = -65603.39999999999
    "H126" = 3251
    "I126" = 3251
    "K126" = 9753
    "M126" = -7283
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
$clears = @("N68", "N71")
foreach ($ref in $clears) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H2" = 247500
    "I2" = 247500
    "K2" = 247500
    "M2" = -247388
    "H62" = 7424.625
    "J62" = 8599.5
    "L62" = 8599.5
    "N62" = -9847.5
    "H65" = 7424.625
    "J65" = 8599.5
    "L65" = 42997.5
    "N65" = -49237.5
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
